$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (pushes existing row 10 "fossil_routes" and below down to row 11+)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new parameter
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true

# Update the sheet dimension-related reference is handled automatically by Excel.
